$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2413793103448276
$ws.Range("C2").Value = 0.3448275862068966
$ws.Range("J2").Value = 0.103448275862069
$ws.Range("P2").Value = 0.1724137931034483
$ws.Range("S2").Value = 0.1379310344827586
$ws.Range("J3").Value = 0.1
$ws.Range("P3").Value = 0.5
$ws.Range("S3").Value = 0.4
$ws.Range("P4").Value = 0.3333333333333333
$ws.Range("S4").Value = 0.6666666666666666
$ws.Range("F6").Value = 0.04166666666666666
$ws.Range("J6").Value = 0.2083333333333333
$ws.Range("Q6").Value = 0.1666666666666667
$ws.Range("R6").Value = 0.04166666666666666
$ws.Range("S6").Value = 0.5416666666666666
$ws.Range("B7").Value = 0.125
$ws.Range("J7").Value = 0.04166666666666666
$ws.Range("Q7").Value = 0.04166666666666666
$ws.Range("S7").Value = 0.7916666666666666
$ws.Range("B8").Value = 0.07216494845360824
$ws.Range("D8").Value = 0.01030927835051546
$ws.Range("F8").Value = 0.04123711340206185
$ws.Range("J8").Value = 0.07216494845360824
$ws.Range("O8").Value = 0.02061855670103093
$ws.Range("Q8").Value = 0.06185567010309279
$ws.Range("R8").Value = 0.07216494845360824
$ws.Range("S8").Value = 0.6494845360824743
$ws.Range("B9").Value = 0.09090909090909091
$ws.Range("F9").Value = 0.04545454545454546
$ws.Range("J9").Value = 0.1363636363636364
$ws.Range("O9").Value = 0.09090909090909091
$ws.Range("Q9").Value = 0.1363636363636364
$ws.Range("R9").Value = 0.04545454545454546
$ws.Range("S9").Value = 0.4545454545454545
$ws.Range("B10").Value = 0.07086614173228346
$ws.Range("D10").Value = 0.01574803149606299
$ws.Range("F10").Value = 0.1023622047244094
$ws.Range("J10").Value = 0.1259842519685039
$ws.Range("O10").Value = 0.01574803149606299
$ws.Range("Q10").Value = 0.1102362204724409
$ws.Range("R10").Value = 0.07086614173228346
$ws.Range("S10").Value = 0.4881889763779528
$ws.Range("G11").Value = 0.1851851851851852
$ws.Range("J11").Value = 0.03703703703703703
$ws.Range("K11").Value = 0.1481481481481481
$ws.Range("L11").Value = 0.5925925925925926
$ws.Range("S11").Value = 0.03703703703703703
$ws.Range("G12").Value = 0.8823529411764706
$ws.Range("S12").Value = 0.1176470588235294
$ws.Range("G13").Value = 0.8
$ws.Range("S13").Value = 0.2
$ws.Range("H15").Value = 0.09523809523809523
$ws.Range("I15").Value = 0.1428571428571428
$ws.Range("J15").Value = 0.2857142857142857
$ws.Range("S15").Value = 0.4761904761904762
$ws.Range("H16").Value = 0.5454545454545454
$ws.Range("J16").Value = 0.1818181818181818
$ws.Range("M16").Value = 0.09090909090909091
$ws.Range("O16").Value = 0.09090909090909091
$ws.Range("S16").Value = 0.09090909090909091
$ws.Range("H17").Value = 0.25
$ws.Range("I17").Value = 0.1428571428571428
$ws.Range("J17").Value = 0.3214285714285715
$ws.Range("K17").Value = 0.03571428571428571
$ws.Range("O17").Value = 0.1428571428571428
$ws.Range("S17").Value = 0.1071428571428571
$ws.Range("H18").Value = 0.2222222222222222
$ws.Range("I18").Value = 0.1666666666666667
$ws.Range("J18").Value = 0.2777777777777778
$ws.Range("K18").Value = 0.05555555555555555
$ws.Range("M18").Value = 0.05555555555555555
$ws.Range("O18").Value = 0.05555555555555555
$ws.Range("S18").Value = 0.1666666666666667
$ws.Range("F19").Value = 0.008695652173913044
$ws.Range("H19").Value = 0.3478260869565217
$ws.Range("I19").Value = 0.05652173913043478
$ws.Range("J19").Value = 0.3043478260869565
$ws.Range("K19").Value = 0.08260869565217391
$ws.Range("M19").Value = 0.01739130434782609
$ws.Range("O19").Value = 0.03043478260869565
$ws.Range("S19").Value = 0.1521739130434783
